# Maschien Vision approach to segemntation task
# Update the Model_Prediction (column C) values on Sheet1 to reflect the
# updated model predictions for a number of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    5   = 3
    8   = 1
    10  = 4
    12  = 4
    13  = 4
    16  = 4
    19  = 3
    20  = 4
    28  = 3
    29  = 4
    32  = 4
    34  = 3
    36  = 4
    46  = 3
    49  = 4
    54  = 1
    57  = 4
    58  = 3
    59  = 4
    60  = 3
    61  = 4
    63  = 4
    64  = 1
    66  = 3
    69  = 3
    70  = 3
    77  = 3
    88  = 3
    90  = 1
    92  = 4
    93  = 1
    98  = 4
    114 = 4
    116 = 3
    126 = 1
    134 = 3
    135 = 4
    143 = 4
    144 = 4
    146 = 3
    149 = 2
    152 = 2
    153 = 1
    159 = 3
    161 = 3
    164 = 3
}

foreach ($row in $changes.Keys) {
    $ws.Range("C$row").Value = $changes[$row]
}
